$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update clearCount / chance values per diff
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 4
$ws.Range("B5").Value = 7
$ws.Range("C6").Value = 6

# Update the active cell selection to E8
$ws.Range("E8").Select()
